$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions) updates
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value = 13268
$ws1.Range("F10").Value = 1189
$ws1.Range("F11").Value = 994
$ws1.Range("F12").Value = 13802
$ws1.Range("F13").Value = 14435
$ws1.Range("F22").Value = 1099
$ws1.Range("F25").Value = 5487
$ws1.Range("F27").Value = 573
$ws1.Range("F28").Value = 343
$ws1.Range("F29").Value = 26
$ws1.Range("F30").Value = 92

# Sheet "全部类型" (All types) updates
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value = 13268
$ws4.Range("F11").Value = 1189
$ws4.Range("F12").Value = 994
$ws4.Range("F13").Value = 13802
$ws4.Range("F14").Value = 14435
$ws4.Range("F23").Value = 1099
$ws4.Range("F26").Value = 5487
$ws4.Range("F28").Value = 573
$ws4.Range("F29").Value = 343
$ws4.Range("F30").Value = 26
$ws4.Range("F31").Value = 92
